$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 384.725
$ws.Range("I15").Value = 384.725
$ws.Range("K15").Value = 1154.175
$ws.Range("M15").Value = -985.1750000000002
$ws.Range("H51").Value = 17646.809
$ws.Range("I51").Value = 15232.167
$ws.Range("J51").Value = 18612.666
$ws.Range("K51").Value = 15232.167
$ws.Range("L51").Value = 18612.666
$ws.Range("M51").Value = -14748.167
$ws.Range("N51").Value = -19580.666
$ws.Range("H52").Value = 1444.5
$ws.Range("I52").Value = 999
$ws.Range("J52").Value = 1890
$ws.Range("K52").Value = 2997
$ws.Range("L52").Value = 5670
$ws.Range("M52").Value = -2837
$ws.Range("N52").Value = -5990
$ws.Range("H55").Value = 651.5
$ws.Range("I55").Value = 318.83334
$ws.Range("J55").Value = 1649.5
$ws.Range("K55").Value = 318.83334
$ws.Range("L55").Value = 1649.5
$ws.Range("M55").Value = -104.83334
$ws.Range("N55").Value = -2077.5
$ws.Range("H64").Value = 2500
$ws.Range("I64").Value = 2500
$ws.Range("K64").Value = 2500
$ws.Range("M64").Value = -2252
$ws.Range("H67").Value = 2500
$ws.Range("I67").Value = 2500
$ws.Range("K67").Value = 2500
$ws.Range("M67").Value = -1642
$ws.Range("H70").Value = 5649.0356
$ws.Range("I70").Value = 6078
$ws.Range("J70").Value = 5154.077
$ws.Range("K70").Value = 18234
$ws.Range("L70").Value = 15462.231
$ws.Range("M70").Value = -17964
$ws.Range("N70").Value = -16002.231
$ws.Range("H73").Value = 5649.0356
$ws.Range("I73").Value = 6078
$ws.Range("J73").Value = 5154.077
$ws.Range("K73").Value = 18234
$ws.Range("L73").Value = 15462.231
$ws.Range("M73").Value = -17298
$ws.Range("N73").Value = -17334.231
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H106").Value = 45486590
$ws.Range("I106").Value = 55586720
$ws.Range("J106").Value = 36000
$ws.Range("K106").Value = 55586720
$ws.Range("L106").Value = 36000
$ws.Range("M106").Value = -55586089
$ws.Range("N106").Value = -37262
$ws.Range("H113").Value = 4000
$ws.Range("I113").Value = 4000
$ws.Range("K113").Value = 4000
$ws.Range("M113").Value = -746
$ws.Range("H116").Value = 3435.3333
$ws.Range("I116").Value = 2650
$ws.Range("K116").Value = 2650
$ws.Range("M116").Value = 792
$ws.Range("H125").Value = 3714
$ws.Range("I125").Value = 3999.75
$ws.Range("J125").Value = 3333
$ws.Range("K125").Value = 35997.75
$ws.Range("L125").Value = 29997
$ws.Range("M125").Value = -33537.75
$ws.Range("N125").Value = -34917
$ws.Range("H131").Value = 2328.7368
$ws.Range("I131").Value = 283.06668
$ws.Range("K131").Value = 849.2000400000001
$ws.Range("M131").Value = 4190.79996
$ws.Range("H132").Value = 4058.9412
$ws.Range("I132").Value = 1427.5454
$ws.Range("J132").Value = 8883.166999999999
$ws.Range("K132").Value = 4282.6362
$ws.Range("L132").Value = 26649.501
$ws.Range("M132").Value = -1752.6362
$ws.Range("N132").Value = -31709.501
$ws.Range("H135").Value = 2529.7144
$ws.Range("I135").Value = 1451.4166
$ws.Range("K135").Value = 13062.7494
$ws.Range("M135").Value = -10527.7494
$ws.Range("H137").Value = 4830.615
$ws.Range("I137").Value = 8200
$ws.Range("K137").Value = 24600
$ws.Range("M137").Value = -22050
$ws.Range("H138").Value = 3103.484
$ws.Range("J138").Value = 5580.8335
$ws.Range("L138").Value = 16742.5005
$ws.Range("N138").Value = -27022.5005
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 6524.75
$ws.Range("J10").Value = 3674.75
$ws.Range("L10").Value = 3674.75
$ws.Range("N10").Value = -4014.75
$ws.Range("H24").Value = 45000
$ws.Range("J24").Value = 45000
$ws.Range("L24").Value = 45000
$ws.Range("N24").Value = -45748
$ws.Range("H32").Value = 372033.4
$ws.Range("I32").Value = 1350.125
$ws.Range("K32").Value = 1350.125
$ws.Range("M32").Value = -1063.125
$ws.Range("H45").Value = 1917.2858
$ws.Range("I45").Value = 1917.2858
$ws.Range("K45").Value = 1917.2858
$ws.Range("M45").Value = -1540.2858
$ws.Range("H63").Value = 8331.166999999999
$ws.Range("I63").Value = 5975
$ws.Range("K63").Value = 5975
$ws.Range("M63").Value = -5289
$ws.Range("H66").Value = 8331.166999999999
$ws.Range("I66").Value = 5975
$ws.Range("K66").Value = 29875
$ws.Range("M66").Value = -26443
$ws.Range("H74").Value = 1462.2354
$ws.Range("I74").Value = 1138.4828
$ws.Range("K74").Value = 1138.4828
$ws.Range("M74").Value = -264.4828
$ws.Range("H77").Value = 1462.2354
$ws.Range("I77").Value = 1138.4828
$ws.Range("K77").Value = 5692.414
$ws.Range("M77").Value = -1324.414
$ws.Range("H88").Value = 3331.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3331.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3331.5
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -4143.5
$ws.Range("H91").Value = 3331.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3331.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3331.5
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -6139.5
$ws.Range("H95").Value = 46666.668
$ws.Range("J95").Value = 46666.668
$ws.Range("L95").Value = 46666.668
$ws.Range("N95").Value = -52158.668
$ws.Range("H97").Value = 532.5833
$ws.Range("I97").Value = 577.3333
$ws.Range("K97").Value = 577.3333
$ws.Range("M97").Value = -81.33330000000001
$ws.Range("H100").Value = 45000
$ws.Range("J100").Value = 45000
$ws.Range("L100").Value = 45000
$ws.Range("N100").Value = -47164
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3254.0625
$ws.Range("I20").Value = 2589.6924
$ws.Range("J20").Value = 6133
$ws.Range("K20").Value = 2589.6924
$ws.Range("L20").Value = 6133
$ws.Range("M20").Value = -2342.6924
$ws.Range("N20").Value = -6627
$ws.Range("H64").Value = 1170
$ws.Range("I64").Value = 1187
$ws.Range("J64").Value = 1158.6666
$ws.Range("K64").Value = 1187
$ws.Range("L64").Value = 1158.6666
$ws.Range("M64").Value = -962
$ws.Range("N64").Value = -1608.6666
$ws.Range("H67").Value = 1170
$ws.Range("I67").Value = 1187
$ws.Range("J67").Value = 1158.6666
$ws.Range("K67").Value = 1187
$ws.Range("L67").Value = 1158.6666
$ws.Range("M67").Value = -407
$ws.Range("N67").Value = -2718.6666
$ws.Range("H82").Value = 221279.8
$ws.Range("J82").Value = 536699.5
$ws.Range("L82").Value = 536699.5
$ws.Range("N82").Value = -537465.5
$ws.Range("H85").Value = 221279.8
$ws.Range("J85").Value = 536699.5
$ws.Range("L85").Value = 536699.5
$ws.Range("N85").Value = -539351.5
$ws.Range("H86").Value = 3000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 3000
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 3000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 15000
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -26232
$ws.Range("H94").Value = 1081.1428
$ws.Range("I94").Value = 1081.1428
$ws.Range("K94").Value = 1081.1428
$ws.Range("M94").Value = -630.1428000000001
$ws.Range("H107").Value = 1292.7222
$ws.Range("I107").Value = 1292.7222
$ws.Range("K107").Value = 1292.7222
$ws.Range("M107").Value = 627.2778000000001
$ws.Range("H134").Value = 4434.2
$ws.Range("I134").Value = 4160.875
$ws.Range("K134").Value = 12482.625
$ws.Range("M134").Value = -9947.625
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 850
$ws.Range("I4").Value = 850
$ws.Range("K4").Value = 850
$ws.Range("M4").Value = -738
$ws.Range("H16").Value = 22729052
$ws.Range("I16").Value = 29413304
$ws.Range("J16").Value = 2593
$ws.Range("K16").Value = 29413304
$ws.Range("L16").Value = 2593
$ws.Range("M16").Value = -29413017
$ws.Range("N16").Value = -3167
$ws.Range("H31").Value = 4783.4
$ws.Range("I31").Value = 4915
$ws.Range("J31").Value = 4717.6
$ws.Range("K31").Value = 4915
$ws.Range("L31").Value = 4717.6
$ws.Range("M31").Value = -4620
$ws.Range("N31").Value = -5307.6
$ws.Range("H34").Value = 4783.4
$ws.Range("I34").Value = 4915
$ws.Range("J34").Value = 4717.6
$ws.Range("K34").Value = 4915
$ws.Range("L34").Value = 4717.6
$ws.Range("M34").Value = -4713
$ws.Range("N34").Value = -5121.6
$ws.Range("H62").Value = 72501.836
$ws.Range("I62").Value = 7449
$ws.Range("J62").Value = 105028.25
$ws.Range("K62").Value = 7449
$ws.Range("L62").Value = 105028.25
$ws.Range("M62").Value = -6825
$ws.Range("N62").Value = -106276.25
$ws.Range("H65").Value = 72501.836
$ws.Range("I65").Value = 7449
$ws.Range("J65").Value = 105028.25
$ws.Range("K65").Value = 37245
$ws.Range("L65").Value = 525141.25
$ws.Range("M65").Value = -34125
$ws.Range("N65").Value = -531381.25
$ws.Range("H95").Value = 26787.285
$ws.Range("J95").Value = 26787.285
$ws.Range("L95").Value = 26787.285
$ws.Range("N95").Value = -32279.285
$ws.Range("H105").Value = 1831.6666
$ws.Range("I105").Value = 1500
$ws.Range("J105").Value = 1997.5
$ws.Range("K105").Value = 1500
$ws.Range("L105").Value = 1997.5
$ws.Range("M105").Value = 247
$ws.Range("N105").Value = -5491.5
$ws.Range("H113").Value = 22729052
$ws.Range("I113").Value = 29413304
$ws.Range("J113").Value = 2593
$ws.Range("K113").Value = 29413304
$ws.Range("L113").Value = 2593
$ws.Range("M113").Value = -29411134
$ws.Range("N113").Value = -6933
$ws.Range("H134").Value = 4158.909
$ws.Range("I134").Value = 4678.4287
$ws.Range("J134").Value = 3249.75
$ws.Range("K134").Value = 14035.2861
$ws.Range("L134").Value = 9749.25
$ws.Range("M134").Value = -11500.2861
$ws.Range("N134").Value = -14819.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 68
$ws.Range("I2").Value = 54.8
$ws.Range("J2").Value = 90
$ws.Range("K2").Value = 328.8
$ws.Range("L2").Value = 540
$ws.Range("M2").Value = -215.8
$ws.Range("N2").Value = -766
$ws.Range("H5").Value = 519.4
$ws.Range("J5").Value = 482.33334
$ws.Range("L5").Value = 1447.00002
$ws.Range("N5").Value = -1671.00002
$ws.Range("H14").Value = 5000
$ws.Range("I14").Value = 5000
$ws.Range("K14").Value = 15000
$ws.Range("M14").Value = -14827
$ws.Range("H38").Value = 147.15384
$ws.Range("I38").Value = 82.14286
$ws.Range("K38").Value = 246.42858
$ws.Range("M38").Value = 100.57142
$ws.Range("H62").Value = 4980.769
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4980.769
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 14942.307
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -16314.307
$ws.Range("H65").Value = 4980.769
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4980.769
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 44826.921
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -51690.921
$ws.Range("H68").Value = 1693.2609
$ws.Range("J68").Value = 1811.5333
$ws.Range("L68").Value = 5434.5999
$ws.Range("N68").Value = -7056.5999
$ws.Range("H71").Value = 1693.2609
$ws.Range("J71").Value = 1811.5333
$ws.Range("L71").Value = 16303.7997
$ws.Range("N71").Value = -24415.7997
$ws.Range("H122").Value = 1532.625
$ws.Range("I122").Value = 1499
$ws.Range("J122").Value = 1543.8334
$ws.Range("K122").Value = 13491
$ws.Range("L122").Value = 13894.5006
$ws.Range("M122").Value = -11041
$ws.Range("N122").Value = -18794.5006
$ws.Range("H131").Value = 1438.9799
$ws.Range("J131").Value = 1445.6598
$ws.Range("L131").Value = 4336.9794
$ws.Range("N131").Value = -14416.9794
$ws.Range("H132").Value = 4228.4614
$ws.Range("I132").Value = 2781.5
$ws.Range("J132").Value = 5468.7144
$ws.Range("K132").Value = 25033.5
$ws.Range("L132").Value = 49218.4296
$ws.Range("M132").Value = -22503.5
$ws.Range("N132").Value = -54278.4296
$ws.Range("H135").Value = 519.4
$ws.Range("J135").Value = 482.33334
$ws.Range("L135").Value = 4341.00006
$ws.Range("N135").Value = -9411.00006
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 46001.332
$ws.Range("J5").Value = 49002
$ws.Range("L5").Value = 49002
$ws.Range("N5").Value = -49226
$ws.Range("H17").Value = 6981.25
$ws.Range("J17").Value = 3454.5
$ws.Range("L17").Value = 3454.5
$ws.Range("N17").Value = -3790.5
$ws.Range("H70").Value = 2659.8
$ws.Range("I70").Value = 2699.75
$ws.Range("J70").Value = 2500
$ws.Range("K70").Value = 2699.75
$ws.Range("L70").Value = 2500
$ws.Range("M70").Value = -2429.75
$ws.Range("N70").Value = -3040
$ws.Range("H73").Value = 2659.8
$ws.Range("I73").Value = 2699.75
$ws.Range("J73").Value = 2500
$ws.Range("K73").Value = 2699.75
$ws.Range("L73").Value = 2500
$ws.Range("M73").Value = -1763.75
$ws.Range("N73").Value = -4372
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -33744
$ws.Range("H97").Value = 5785.6665
$ws.Range("I97").Value = 1497
$ws.Range("J97").Value = 6643.4
$ws.Range("K97").Value = 1497
$ws.Range("L97").Value = 6643.4
$ws.Range("M97").Value = -1001
$ws.Range("N97").Value = -7635.4
$ws.Range("H113").Value = 3013.625
$ws.Range("I113").Value = 1652.25
$ws.Range("K113").Value = 1652.25
$ws.Range("M113").Value = 517.75
$ws.Range("H122").Value = 64935.625
$ws.Range("I122").Value = 2598.0667
$ws.Range("K122").Value = 7794.2001
$ws.Range("M122").Value = -5344.2001
$ws.Range("H126").Value = 3717.926
$ws.Range("I126").Value = 3220.5557
$ws.Range("K126").Value = 9661.667099999999
$ws.Range("M126").Value = -7191.667099999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 3502.6667
$ws.Range("I17").Value = 2254
$ws.Range("J17").Value = 6000
$ws.Range("K17").Value = 2254
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = -2084
$ws.Range("N17").Value = -6340
$ws.Range("H22").Value = 1072.6
$ws.Range("I22").Value = 799.75
$ws.Range("K22").Value = 799.75
$ws.Range("M22").Value = -504.75
$ws.Range("H27").Value = 1072.6
$ws.Range("I27").Value = 799.75
$ws.Range("K27").Value = 799.75
$ws.Range("M27").Value = -692.75
$ws.Range("H39").Value = 43500
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 43500
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 43500
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -44420
$ws.Range("H46").Value = 3114.1428
$ws.Range("I46").Value = 733
$ws.Range("K46").Value = 733
$ws.Range("M46").Value = -545
$ws.Range("H61").Value = 12347738
$ws.Range("I61").Value = 12347738
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 12347738
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -12347536
$ws.Range("N61").ClearContents()
$ws.Range("H82").Value = 50860.695
$ws.Range("J82").Value = 103705.4
$ws.Range("L82").Value = 103705.4
$ws.Range("N82").Value = -104427.4
$ws.Range("H85").Value = 50860.695
$ws.Range("J85").Value = 103705.4
$ws.Range("L85").Value = 103705.4
$ws.Range("N85").Value = -106201.4
$ws.Range("H113").Value = 12347738
$ws.Range("I113").Value = 12347738
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 12347738
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -12345568
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 5000
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 3895.1667
$ws.Range("I136").Value = 3895.1667
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 11685.5001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -9135.500100000001
$ws.Range("N136").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 4996
$ws.Range("I9").Value = 4996
$ws.Range("K9").Value = 4996
$ws.Range("M9").Value = -4856
$ws.Range("H26").Value = 413799.6
$ws.Range("I26").Value = 14499.5
$ws.Range("J26").Value = 679999.7
$ws.Range("K26").Value = 14499.5
$ws.Range("L26").Value = 679999.7
$ws.Range("M26").Value = -14206.5
$ws.Range("N26").Value = -680585.7
$ws.Range("H61").Value = 39000
$ws.Range("I61").Value = 39000
$ws.Range("K61").Value = 39000
$ws.Range("M61").Value = -38708
$ws.Range("H62").Value = 5000.2
$ws.Range("I62").Value = 2000.4
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 2000.4
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -1376.4
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 5000.2
$ws.Range("I65").Value = 2000.4
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 10002
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -6882
$ws.Range("N65").Value = -46240
$ws.Range("H100").Value = 814.4286
$ws.Range("I100").Value = 783.5
$ws.Range("K100").Value = 1567
$ws.Range("M100").Value = -1026
$ws.Range("H122").Value = 4066.6365
$ws.Range("I122").Value = 4539.857
$ws.Range("K122").Value = 13619.571
$ws.Range("M122").Value = -11169.571
$ws.Range("H126").Value = 447.5
$ws.Range("I126").Value = 447.5
$ws.Range("K126").Value = 1342.5
$ws.Range("M126").Value = 1127.5
$ws.Range("H136").Value = 5254.96
$ws.Range("I136").Value = 4335.1816
$ws.Range("K136").Value = 13005.5448
$ws.Range("M136").Value = -10455.5448
